$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 256.25
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 200
$ws.Range("M33").Value = 29

$ws.Range("H51").Value = 2995.5
$ws.Range("I51").Value = 2995.5
$ws.Range("K51").Value = 2995.5
$ws.Range("M51").Value = -2511.5

$ws.Range("H92").Value = 1409.375
$ws.Range("I92").Value = 1490.0667
$ws.Range("K92").Value = 1490.0667
$ws.Range("M92").Value = -242.0667000000001

$ws.Range("H96").Value = 1809.5416
$ws.Range("I96").Value = 786.75
$ws.Range("J96").Value = 2832.3333
$ws.Range("K96").Value = 2360.25
$ws.Range("L96").Value = 8496.999899999999
$ws.Range("M96").Value = -987.25
$ws.Range("N96").Value = -11242.9999

$ws.Range("H100").Value = 3544.4856
$ws.Range("I100").Value = 1595.125
$ws.Range("J100").Value = 5186.0527
$ws.Range("K100").Value = 1595.125
$ws.Range("L100").Value = 5186.0527
$ws.Range("M100").Value = -1054.125
$ws.Range("N100").Value = -6268.0527

$ws.Range("H107").Value = 11999.8
$ws.Range("I107").Value = 11999.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 11999.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -10079.8
$ws.Range("N107").ClearContents()

$ws.Range("H116").Value = 15900.923
$ws.Range("I116").Value = 19620
$ws.Range("K116").Value = 19620
$ws.Range("M116").Value = -16178

$ws.Range("H132").Value = 8516.538
$ws.Range("I132").Value = 9830.227999999999
$ws.Range("K132").Value = 29490.684
$ws.Range("M132").Value = -26960.684

$ws.Range("H135").Value = 5112.385
$ws.Range("I135").Value = 832.86664
$ws.Range("J135").Value = 10948.091
$ws.Range("K135").Value = 7495.79976
$ws.Range("L135").Value = 98532.819
$ws.Range("M135").Value = -4960.79976
$ws.Range("N135").Value = -103602.819

$ws.Range("H137").Value = 6049.7915
$ws.Range("I137").Value = 1914.1428
$ws.Range("J137").Value = 34999.332
$ws.Range("K137").Value = 5742.428400000001
$ws.Range("L137").Value = 104997.996
$ws.Range("M137").Value = -3192.428400000001
$ws.Range("N137").Value = -110097.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 183566.16
$ws.Range("I32").Value = 195382.98
$ws.Range("K32").Value = 195382.98
$ws.Range("M32").Value = -195095.98

$ws.Range("H61").Value = 20081.834
$ws.Range("I61").Value = 23099.2
$ws.Range("K61").Value = 23099.2
$ws.Range("M61").Value = -22887.2

$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344

$ws.Range("H122").Value = 2647.8572
$ws.Range("I122").Value = 2212.8572
$ws.Range("J122").Value = 3082.8572
$ws.Range("K122").Value = 6638.571599999999
$ws.Range("L122").Value = 9248.571599999999
$ws.Range("M122").Value = -4188.571599999999
$ws.Range("N122").Value = -14148.5716

$ws.Range("H136").Value = 20081.834
$ws.Range("I136").Value = 23099.2
$ws.Range("K136").Value = 69297.60000000001
$ws.Range("M136").Value = -66747.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1413.4722
$ws.Range("I107").Value = 1192.2609
$ws.Range("K107").Value = 1192.2609
$ws.Range("M107").Value = 727.7391

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3130.2222
$ws.Range("I31").Value = 2896.75
$ws.Range("K31").Value = 2896.75
$ws.Range("M31").Value = -2601.75

$ws.Range("H34").Value = 3130.2222
$ws.Range("I34").Value = 2896.75
$ws.Range("K34").Value = 2896.75
$ws.Range("M34").Value = -2694.75

$ws.Range("H52").Value = 73180
$ws.Range("J52").Value = 84770
$ws.Range("L52").Value = 84770
$ws.Range("N52").Value = -85358

$ws.Range("H58").Value = 7524.357
$ws.Range("I58").Value = 6575.75
$ws.Range("J58").Value = 7903.8
$ws.Range("K58").Value = 6575.75
$ws.Range("L58").Value = 7903.8
$ws.Range("M58").Value = -6372.75
$ws.Range("N58").Value = -8309.799999999999

$ws.Range("H99").Value = 28694.75
$ws.Range("J99").Value = 3112.25
$ws.Range("L99").Value = 3112.25
$ws.Range("N99").Value = -6108.25

$ws.Range("H122").Value = 27211.2
$ws.Range("I122").Value = 2179.6
$ws.Range("J122").Value = 52242.8
$ws.Range("K122").Value = 6538.799999999999
$ws.Range("L122").Value = 156728.4
$ws.Range("M122").Value = -4088.799999999999
$ws.Range("N122").Value = -161628.4

$ws.Range("H126").Value = 28694.75
$ws.Range("J126").Value = 3112.25
$ws.Range("L126").Value = 9336.75
$ws.Range("N126").Value = -14276.75

$ws.Range("H132").Value = 10813.272
$ws.Range("I132").Value = 10813.272
$ws.Range("K132").Value = 32439.816
$ws.Range("M132").Value = -29909.816

$ws.Range("H134").Value = 3348.25
$ws.Range("I134").Value = 3406.5
$ws.Range("K134").Value = 10219.5
$ws.Range("M134").Value = -7684.5

$ws.Range("H136").Value = 7524.357
$ws.Range("I136").Value = 6575.75
$ws.Range("J136").Value = 7903.8
$ws.Range("K136").Value = 19727.25
$ws.Range("L136").Value = 23711.4
$ws.Range("M136").Value = -17177.25
$ws.Range("N136").Value = -28811.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6902.4443
$ws.Range("I3").Value = 2749.4546
$ws.Range("J3").Value = 13428.571
$ws.Range("K3").Value = 8248.363799999999
$ws.Range("L3").Value = 40285.713
$ws.Range("M3").Value = -8136.363799999999
$ws.Range("N3").Value = -40509.713

$ws.Range("H26").Value = 356.63635
$ws.Range("J26").Value = 1101
$ws.Range("L26").Value = 3303
$ws.Range("N26").Value = -3879

$ws.Range("H81").Value = 1667.8889
$ws.Range("I81").Value = 1003.6667
$ws.Range("K81").Value = 3011.0001
$ws.Range("M81").Value = -1888.0001

$ws.Range("H84").Value = 1667.8889
$ws.Range("I84").Value = 1003.6667
$ws.Range("K84").Value = 9033.0003
$ws.Range("M84").Value = -3417.0003

$ws.Range("H131").Value = 3423.3333
$ws.Range("J131").Value = 3687.2683
$ws.Range("L131").Value = 11061.8049
$ws.Range("N131").Value = -21141.8049

$ws.Range("H132").Value = 5572.3076
$ws.Range("J132").Value = 6736.1577
$ws.Range("L132").Value = 60625.41929999999
$ws.Range("N132").Value = -65685.41929999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12892.652
$ws.Range("I70").Value = 13751.723
$ws.Range("K70").Value = 13751.723
$ws.Range("M70").Value = -13481.723

$ws.Range("H73").Value = 12892.652
$ws.Range("I73").Value = 13751.723
$ws.Range("K73").Value = 13751.723
$ws.Range("M73").Value = -12815.723

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H97").Value = 2083.5264
$ws.Range("I97").Value = 1952.0625
$ws.Range("K97").Value = 1952.0625
$ws.Range("M97").Value = -1456.0625

$ws.Range("H102").Value = 2493
$ws.Range("I102").Value = 992.06665
$ws.Range("K102").Value = 992.06665
$ws.Range("M102").Value = 629.93335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1070.4286
$ws.Range("I16").Value = 1070.4286
$ws.Range("K16").Value = 1070.4286
$ws.Range("M16").Value = -900.4286

$ws.Range("H40").Value = 5572.273
$ws.Range("I40").Value = 3150
$ws.Range("J40").Value = 6956.4287
$ws.Range("K40").Value = 3150
$ws.Range("L40").Value = 6956.4287
$ws.Range("M40").Value = -3014
$ws.Range("N40").Value = -7228.4287

$ws.Range("H46").Value = 2934.1353
$ws.Range("J46").Value = 3633.1924
$ws.Range("L46").Value = 3633.1924
$ws.Range("N46").Value = -4009.1924

$ws.Range("H82").Value = 11508.286
$ws.Range("I82").Value = 13010.5
$ws.Range("K82").Value = 13010.5
$ws.Range("M82").Value = -12649.5

$ws.Range("H85").Value = 11508.286
$ws.Range("I85").Value = 13010.5
$ws.Range("K85").Value = 13010.5
$ws.Range("M85").Value = -11762.5

$ws.Range("H93").Value = 7659
$ws.Range("I93").Value = 7032.4
$ws.Range("K93").Value = 7032.4
$ws.Range("M93").Value = -5784.4

$ws.Range("H122").Value = 7095
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws.Range("H132").Value = 4665.6665
$ws.Range("I132").Value = 4665.6665
$ws.Range("K132").Value = 13996.9995
$ws.Range("M132").Value = -11466.9995

$ws.Range("H136").Value = 6275.7
$ws.Range("I136").Value = 3623.6667
$ws.Range("K136").Value = 10871.0001
$ws.Range("M136").Value = -8321.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 30000
$ws.Range("I2").Value = 30000
$ws.Range("K2").Value = 30000
$ws.Range("M2").Value = -29888

$ws.Range("H122").Value = 73268.5
$ws.Range("I122").Value = 2392.3333
$ws.Range("J122").Value = 115794.2
$ws.Range("K122").Value = 7176.999899999999
$ws.Range("L122").Value = 347382.6
$ws.Range("M122").Value = -4726.999899999999
$ws.Range("N122").Value = -352282.6

$ws.Range("H132").Value = 1913.0294
$ws.Range("I132").Value = 1755.5161
$ws.Range("J132").Value = 3540.6667
$ws.Range("K132").Value = 5266.5483
$ws.Range("L132").Value = 10622.0001
$ws.Range("M132").Value = -2736.5483
$ws.Range("N132").Value = -15682.0001
